$wb = $excel.ActiveWorkbook

# Sheet "展览" updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 56
$ws1.Range("F4").Value = 3584
$ws1.Range("F6").Value = 433
$ws1.Range("F7").Value = 2
$ws1.Range("F11").Value = 1324
$ws1.Range("F13").Value = 1896

# Sheet "全部类型" updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 56
$ws4.Range("F4").Value = 3584
$ws4.Range("F6").Value = 433
$ws4.Range("F7").Value = 2
$ws4.Range("F14").Value = 1324
$ws4.Range("F16").Value = 1896
